$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("levels")
$ws.Range("A1").Value = "test"
